$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A179").Value = "Python Developer with GoLang"
$ws.Range("B179").Value = "https://www.dice.com/job-detail/06ed2f6d-b4fc-48f0-b23e-99da261889e5?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C179").Value = "Chicago, Illinois"
$ws.Range("D179").Value = "Full-time, Third Party"
$ws.Range("E179").Value = "Depends on Experience"
$ws.Range("F179").Value = "Dahl Consulting"

$ws.Range("A180").Value = "Golang Developer"
$ws.Range("B180").Value = "https://www.dice.com/job-detail/4bb243d1-ed24-4c99-9fb0-76b8e53da475?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C180").Value = "Hybrid in Dallas, Texas"
$ws.Range("D180").Value = "Contract, Third Party"
$ws.Range("E180").Value = "50 - 55"
$ws.Range("F180").Value = "Stellar Professionals LLC"
